$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D values that look like plain decimal numbers remain stored as text,
# matching the source data (which stores all Price values as text).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '43.361.34'
$ws.Range('E2').Value = '  +2.69%  '
$ws.Range('D3').Value = '2.306.74'
$ws.Range('D4').Value = '0.999'
$ws.Range('D5').Value = '310.44'
$ws.Range('E5').Value = '  +1.25%  '
$ws.Range('D6').Value = '103.04'
$ws.Range('E6').Value = '  +5.56%  '
$ws.Range('D7').Value = '0.534'
$ws.Range('E7').Value = '  +1.49%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('E9').Value = '  +8.36%  '
$ws.Range('D10').Value = '35.65'
$ws.Range('E10').Value = '  +0.63%  '
$ws.Range('E11').Value = '  +2.89%  '
$ws.Range('E12').Value = '  -1.10%  '
$ws.Range('D13').Value = '7.02'
$ws.Range('E13').Value = '  +2.15%  '
$ws.Range('D14').Value = '2.663.12'
$ws.Range('E14').Value = '  +1.56%  '
$ws.Range('E15').Value = '  +1.75%  '
$ws.Range('D16').Value = '2.280.81'
$ws.Range('E16').Value = '  +1.79%  '
$ws.Range('D17').Value = '0.807'
$ws.Range('E17').Value = '  +2.01%  '
$ws.Range('D18').Value = '43.256.05'
$ws.Range('E18').Value = '  +2.73%  '
$ws.Range('E19').Value = '  -0.56%  '
$ws.Range('E20').Value = '  +3.20%  '
$ws.Range('E21').Value = '  +2.76%  '
$ws.Range('D22').Value = '68.11'
$ws.Range('E22').Value = '  +0.52%  '
$ws.Range('D23').Value = '241.76'
$ws.Range('E23').Value = '  +1.94%  '
$ws.Range('E24').Value = '  +0.51%  '
$ws.Range('D25').Value = '2.61'
$ws.Range('E25').Value = '  +0.74%  '
$ws.Range('E26').Value = '  +0.07%  '
$ws.Range('D27').Value = '25.00'
$ws.Range('E27').Value = '  +6.14%  '
$ws.Range('E28').Value = '  +7.89%  '
$ws.Range('D29').Value = '36.61'
$ws.Range('E29').Value = '  -2.10%  '
$ws.Range('D30').Value = '9.66'
$ws.Range('E30').Value = '  +0.65%  '
$ws.Range('D31').Value = '171.83'
$ws.Range('E31').Value = '  +5.61%  '
$ws.Range('D32').Value = '5.28'
$ws.Range('E32').Value = '  +0.21%  '
$ws.Range('E33').Value = '  +0.00%  '
$ws.Range('D34').Value = '2.57'
$ws.Range('E34').Value = '  +7.79%  '
$ws.Range('E35').Value = '  +0.43%  '
$ws.Range('D36').Value = '0.0741'
$ws.Range('E36').Value = '  +0.78%  '
$ws.Range('E37').Value = '  -2.55%  '
$ws.Range('E38').Value = '  +2.88%  '
$ws.Range('E39').Value = '  +1.38%  '
$ws.Range('E40').Value = '  +1.51%  '
$ws.Range('D41').Value = '4.33'
$ws.Range('E41').Value = '  +5.36%  '
$ws.Range('E42').Value = '  -1.59%  '
$ws.Range('D43').Value = '0.0293'
$ws.Range('E43').Value = '  +4.35%  '
$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = '1.971.51'
$ws.Range('E44').Value = '  +1.14%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = '19.14'
$ws.Range('E45').Value = '  +0.24%  '
$ws.Range('D46').Value = '2.99'
$ws.Range('E46').Value = '  +2.45%  '
$ws.Range('D47').Value = '9.97'
$ws.Range('E47').Value = '  -0.01%  '
$ws.Range('D48').Value = '55.45'
$ws.Range('E48').Value = '  +3.22%  '
$ws.Range('B49').Value = 'HuobiToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D49').Value = '2.93'
$ws.Range('E49').Value = '  +1.94%  '
$ws.Range('B50').Value = 'Stacks'
$ws.Range('C50').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D50').Value = '1.58'
$ws.Range('E50').Value = '  +6.86%  '
$ws.Range('D51').Value = '2.531.27'
$ws.Range('E51').Value = '  +1.55%  '
